$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228
$ws.Range("H9").Value = 1082601.9
$ws.Range("I9").Value = 2164767.8
$ws.Range("J9").Value = 435.83334
$ws.Range("K9").Value = 2164767.8
$ws.Range("L9").Value = 435.83334
$ws.Range("M9").Value = -2164598.8
$ws.Range("N9").Value = -773.83334
$ws.Range("H12").Value = 198
$ws.Range("J12").Value = 196.66667
$ws.Range("L12").Value = 196.66667
$ws.Range("N12").Value = -536.6666700000001
$ws.Range("H28").Value = 1397.6364
$ws.Range("J28").Value = 3376.5
$ws.Range("L28").Value = 3376.5
$ws.Range("N28").Value = -4346.5
$ws.Range("H40").Value = 11261.808
$ws.Range("I40").Value = 3933.3
$ws.Range("J40").Value = 15842.125
$ws.Range("K40").Value = 3933.3
$ws.Range("L40").Value = 15842.125
$ws.Range("M40").Value = -3758.3
$ws.Range("N40").Value = -16192.125
$ws.Range("H51").Value = 27826.092
$ws.Range("I51").Value = 9900
$ws.Range("J51").Value = 34548.375
$ws.Range("K51").Value = 9900
$ws.Range("L51").Value = 34548.375
$ws.Range("M51").Value = -9416
$ws.Range("N51").Value = -35516.375
$ws.Range("H80").Value = 1313.5714
$ws.Range("I80").Value = 598.3333
$ws.Range("J80").Value = 1850
$ws.Range("K80").Value = 1794.9999
$ws.Range("L80").Value = 5550
$ws.Range("M80").Value = -796.9999
$ws.Range("N80").Value = -7546
$ws.Range("H83").Value = 1313.5714
$ws.Range("I83").Value = 598.3333
$ws.Range("J83").Value = 1850
$ws.Range("K83").Value = 5384.9997
$ws.Range("L83").Value = 16650
$ws.Range("M83").Value = -392.9997000000003
$ws.Range("N83").Value = -26634
$ws.Range("H88").Value = 2788.2173
$ws.Range("I88").Value = 1694.6666
$ws.Range("J88").Value = 3174.1765
$ws.Range("K88").Value = 1694.6666
$ws.Range("L88").Value = 3174.1765
$ws.Range("M88").Value = -1288.6666
$ws.Range("N88").Value = -3986.1765
$ws.Range("H91").Value = 2788.2173
$ws.Range("I91").Value = 1694.6666
$ws.Range("J91").Value = 3174.1765
$ws.Range("K91").Value = 1694.6666
$ws.Range("L91").Value = 3174.1765
$ws.Range("M91").Value = -290.6666
$ws.Range("N91").Value = -5982.1765
$ws.Range("H96").Value = 1352.7142
$ws.Range("I96").Value = 813.3333
$ws.Range("K96").Value = 2439.9999
$ws.Range("M96").Value = -1066.9999
$ws.Range("H97").Value = 3766.6667
$ws.Range("J97").Value = 3766.6667
$ws.Range("L97").Value = 11300.0001
$ws.Range("N97").Value = -12292.0001
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H125").Value = 982.6667
$ws.Range("I125").Value = 979.4
$ws.Range("K125").Value = 8814.6
$ws.Range("M125").Value = -6354.6
$ws.Range("H137").Value = 2318.92
$ws.Range("I137").Value = 1561.091
$ws.Range("K137").Value = 4683.272999999999
$ws.Range("M137").Value = -2133.272999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 988.2692
$ws.Range("I110").Value = 924.7917
$ws.Range("K110").Value = 924.7917
$ws.Range("M110").Value = 1120.2083

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3925.0667
$ws.Range("I99").Value = 3721.3076
$ws.Range("K99").Value = 3721.3076
$ws.Range("M99").Value = -2223.3076
$ws.Range("H112").Value = 45000
$ws.Range("J112").Value = 45000
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47954

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 99999
$ws.Range("J8").Value = 99999
$ws.Range("L8").Value = 99999
$ws.Range("N8").Value = -100279
$ws.Range("H16").Value = 11687.267
$ws.Range("I16").Value = 950.125
$ws.Range("K16").Value = 950.125
$ws.Range("M16").Value = -663.125
$ws.Range("H86").Value = 14499336
$ws.Range("I86").Value = 25646530
$ws.Range("J86").Value = 7982.7
$ws.Range("K86").Value = 25646530
$ws.Range("L86").Value = 7982.7
$ws.Range("M86").Value = -25645407
$ws.Range("N86").Value = -10228.7
$ws.Range("H89").Value = 14499336
$ws.Range("I89").Value = 25646530
$ws.Range("J89").Value = 7982.7
$ws.Range("K89").Value = 128232650
$ws.Range("L89").Value = 39913.5
$ws.Range("M89").Value = -128227034
$ws.Range("N89").Value = -51145.5
$ws.Range("H113").Value = 11687.267
$ws.Range("I113").Value = 950.125
$ws.Range("K113").Value = 950.125
$ws.Range("M113").Value = 1219.875
$ws.Range("H134").Value = 3956.75
$ws.Range("I134").Value = 3862.0908
$ws.Range("K134").Value = 11586.2724
$ws.Range("M134").Value = -9051.2724

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 11937.25
$ws.Range("J57").Value = 14749.833
$ws.Range("L57").Value = 44249.499
$ws.Range("N57").Value = -45367.499
$ws.Range("H98").Value = 1360.4375
$ws.Range("I98").Value = 1160.125
$ws.Range("J98").Value = 1560.75
$ws.Range("K98").Value = 3480.375
$ws.Range("L98").Value = 4682.25
$ws.Range("M98").Value = -1982.375
$ws.Range("N98").Value = -7678.25
$ws.Range("H118").Value = 7998.5
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 7998.5
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 23995.5
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -26481.5
$ws.Range("H121").Value = 136960.05
$ws.Range("J121").Value = 501362
$ws.Range("L121").Value = 1504086
$ws.Range("N121").Value = -1506706
$ws.Range("H131").Value = 1746.6666
$ws.Range("J131").Value = 2294.5833
$ws.Range("L131").Value = 6883.749899999999
$ws.Range("N131").Value = -16963.7499

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3957.125
$ws.Range("I80").Value = 1971
$ws.Range("J80").Value = 5943.25
$ws.Range("K80").Value = 1971
$ws.Range("L80").Value = 5943.25
$ws.Range("M80").Value = -973
$ws.Range("N80").Value = -7939.25
$ws.Range("H83").Value = 3957.125
$ws.Range("I83").Value = 1971
$ws.Range("J83").Value = 5943.25
$ws.Range("K83").Value = 9855
$ws.Range("L83").Value = 29716.25
$ws.Range("M83").Value = -4863
$ws.Range("N83").Value = -39700.25
$ws.Range("H97").Value = 1071.5
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 547.61536
$ws.Range("I107").Value = 557.7692
$ws.Range("J107").Value = 537.46155
$ws.Range("K107").Value = 557.7692
$ws.Range("L107").Value = 537.46155
$ws.Range("M107").Value = 1362.2308
$ws.Range("N107").Value = -4377.46155

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5694.706
$ws.Range("I46").Value = 7031.4287
$ws.Range("J46").Value = 4759
$ws.Range("K46").Value = 7031.4287
$ws.Range("L46").Value = 4759
$ws.Range("M46").Value = -6843.4287
$ws.Range("N46").Value = -5135
$ws.Range("H74").Value = 26000
$ws.Range("I74").Value = 26000
$ws.Range("K74").Value = 26000
$ws.Range("M74").Value = -25002
$ws.Range("H77").Value = 26000
$ws.Range("I77").Value = 26000
$ws.Range("K77").Value = 78000
$ws.Range("M77").Value = -73008
$ws.Range("H93").Value = 2010.8
$ws.Range("J93").Value = 7004
$ws.Range("L93").Value = 7004
$ws.Range("N93").Value = -9500
